$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: fill A7:E7 with "f"
$ws.Range("A7:E7").Value = "f"

# Row 8: new scouting entry
$ws.Range("A8").Value = "suday"

# Team Number column holds text-formatted numbers (see B2 = "6897" as text).
# Writing a numeric-looking string via .Value would get auto-coerced to a
# number by Excel's type inference, so build it as a text formula first and
# then flatten it to a literal value via copy / paste-special values. This
# keeps the cell a plain shared-string cell with no extra number formatting.
$cell = $ws.Cells.Item(8, 2)
$cell.Formula = '="6897"'
$cell.Copy()
$cell.PasteSpecial(-4163)

$ws.Range("C8").Value = "cool"
$ws.Range("D8").Value = "test"
$ws.Range("E8").Value = "hello"
